# The source sheet stores Price (D) and Volume(1h) (E) as text (not real
# numbers) even though most look numeric. A leading apostrophe forces
# Excel to keep the literal text instead of auto-converting to a number/
# percentage, matching the original "General"-formatted text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple Price (D) / Volume(1h) (E) updates that don't move rows ---
$ws.Range("D2").Value = "'243.47"
$ws.Range("E2").Value = "'-0.73%"

$ws.Range("D3").Value = "'29.75"
$ws.Range("E3").Value = "'12.11%"

$ws.Range("D4").Value = "'5.121"
$ws.Range("E4").Value = "'0.12%"

$ws.Range("D5").Value = "'0.05673"
$ws.Range("E5").Value = "'1.52%"

$ws.Range("D6").Value = "'6.508"
$ws.Range("E6").Value = "'0.51%"

$ws.Range("D7").Value = "'0.8261"
$ws.Range("E7").Value = "'1.11%"

$ws.Range("D8").Value = "'0.8618"
$ws.Range("E8").Value = "'3.22%"

# --- Rows 9-15: "One" moves to the top of this block (row 9), the rest
#     (WazirX..CoinExToken) shift down by one row. Rewrite Coin (B),
#     Link (C), Price (D) and Volume(1h) (E) for each of these rows. ---

$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01015"
$ws.Range("E9").Value = "'1,589.07%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1336"
$ws.Range("E10").Value = "'0.28%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06925"
$ws.Range("E11").Value = "'-0.95%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02858"
$ws.Range("E12").Value = "'-1.05%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09374"
$ws.Range("E13").Value = "'-0.15%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001509"
$ws.Range("E14").Value = "'-0.23%"

$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04146"
$ws.Range("E15").Value = "'-9.43%"

# --- Remaining Price (D) / Volume(1h) (E) updates (rows unaffected by the
#     row-9..15 reshuffle) ---
$ws.Range("D16").Value = "'0.006025"
$ws.Range("E16").Value = "'-1.88%"

$ws.Range("D17").Value = "'3.519"
$ws.Range("E17").Value = "'-3.55%"

$ws.Range("D18").Value = "'3.009"
$ws.Range("E18").Value = "'-0.92%"

$ws.Range("D19").Value = "'2.220"
$ws.Range("E19").Value = "'1.69%"

$ws.Range("D21").Value = "'0.03286"
$ws.Range("E21").Value = "'4.73%"

$ws.Range("D22").Value = "'0.1295"
$ws.Range("E22").Value = "'-0.34%"

$ws.Range("D23").Value = "'3.609"
$ws.Range("E23").Value = "'-3.51%"

$ws.Range("E24").Value = "'-0.05%"

$ws.Range("D25").Value = "'0.001208"
$ws.Range("E25").Value = "'-3.01%"

$ws.Range("D26").Value = "'0.004455"
$ws.Range("E26").Value = "'-1.27%"

$ws.Range("D27").Value = "'0.0001176"
$ws.Range("E27").Value = "'22.49%"

$ws.Range("D28").Value = "'0.0001399"
$ws.Range("E28").Value = "'0.22%"

$ws.Range("D40").Value = "'0.03706"
$ws.Range("E40").Value = "'1.81%"

$ws.Range("D41").Value = "'0.005769"
$ws.Range("E41").Value = "'-6.48%"

$ws.Range("D42").Value = "'0.1056"
$ws.Range("E42").Value = "'0.40%"

$ws.Range("D43").Value = "'0.002303"
$ws.Range("E43").Value = "'-4.03%"

$ws.Range("D44").Value = "'0.009271"
$ws.Range("E44").Value = "'4.64%"

$ws.Range("D45").Value = "'0.00005079"
$ws.Range("E45").Value = "'-5.01%"

$ws.Range("D46").Value = "'0.00000000747"
$ws.Range("E46").Value = "'-0.34%"

$ws.Range("D47").Value = "'0.1006"
$ws.Range("E47").Value = "'-7.66%"

$ws.Range("D48").Value = "'0.002717"
$ws.Range("E48").Value = "'6.45%"

$ws.Range("D49").Value = "'0.00002093"
$ws.Range("E49").Value = "'-0.34%"

$ws.Range("D50").Value = "'0.0001993"
$ws.Range("E50").Value = "'-0.34%"
